# Update "paises" workbook: refresh case counters and re-sort a block of
# countries that tied on "Casos totales" (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the "last updated" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 22:22"

# 2) Row 4 - Estados Unidos: updated counters.
$ws.Cells.Item(4, 2).Value = 493360
$ws.Cells.Item(4, 3).Value = 24794
$ws.Cells.Item(4, 5).Value = 448247
$ws.Cells.Item(4, 7).Value = 1639
$ws.Cells.Item(4, 8).Value = 18330

# 3) Row 28: updated counters.
$ws.Cells.Item(28, 2).Value = 6314
$ws.Cells.Item(28, 3).Value = 95
$ws.Cells.Item(28, 5).Value = 6169

# 4) Rows 180-189: countries tied at ~11-15 total cases got refreshed and
#    re-sorted by "Casos totales" (column B) descending.
$ws.Cells.Item(180, 1).Value = "Granada"
$ws.Cells.Item(180, 2).Value = 14
$ws.Cells.Item(180, 3).Value = 2
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 14
$ws.Cells.Item(180, 6).Value = 2
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "Curazao"
$ws.Cells.Item(181, 2).Value = 14
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 7
$ws.Cells.Item(181, 5).Value = 6
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 1

$ws.Cells.Item(182, 1).Value = "Botsuana"
$ws.Cells.Item(182, 2).Value = 13
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 12
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 1

$ws.Cells.Item(183, 1).Value = "Zimbabue"
$ws.Cells.Item(183, 2).Value = 13
$ws.Cells.Item(183, 3).Value = 2
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 10
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 3

$ws.Cells.Item(184, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(184, 2).Value = 12
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 1
$ws.Cells.Item(184, 5).Value = 11
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Suazilandia"
$ws.Cells.Item(185, 2).Value = 12
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 7
$ws.Cells.Item(185, 5).Value = 5
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

$ws.Cells.Item(186, 1).Value = "Seychelles"
$ws.Cells.Item(186, 2).Value = 11
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 11
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

$ws.Cells.Item(187, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(187, 2).Value = 11
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 11
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

$ws.Cells.Item(188, 1).Value = "Republica del Chad"
$ws.Cells.Item(188, 2).Value = 11
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 2
$ws.Cells.Item(188, 5).Value = 9
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

$ws.Cells.Item(189, 1).Value = "Groenlandia"
$ws.Cells.Item(189, 2).Value = 11
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 11
$ws.Cells.Item(189, 5).Value = 0
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0
